$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new rows 6-9 are exact duplicates of existing rows 4, 2, 3, 5
# (in that order). Copy/paste preserves the "numberStoredAsText" string
# typing (t="str") of the source cells without introducing a new style.
$ws.Range("A4:K4").Copy()
$ws.Range("A6:K6").PasteSpecial()

$ws.Range("A2:K2").Copy()
$ws.Range("A7:K7").PasteSpecial()

$ws.Range("A3:K3").Copy()
$ws.Range("A8:K8").PasteSpecial()

$ws.Range("A5:K5").Copy()
$ws.Range("A9:K9").PasteSpecial()
